$wb = $excel.ActiveWorkbook

# Use the existing "Czech" sheet as the template for the new "Swiss" sheet,
# since both sheets share the same repeaters/layout structure.
$czech = $wb.Worksheets.Item("Czech")

# Copy Czech to the end of the workbook, then rename the copy to "Swiss"
$czech.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# The Swiss market needs two extra repeater rows (PR1DSCH / PR8ASCH) that
# Czech doesn't have: one right after "PR1DS" (row 18) and one right after
# "PR8AS" (row 20, once the first insert has shifted it down). This keeps
# the original formatting of the existing rows intact.
$swiss.Rows.Item(19).Insert()
$swiss.Rows.Item(21).Insert()

# Seed the formatting of the two freshly inserted rows from their neighbors
# so they pick up the same border/style as the rest of the list.
$swiss.Range("A18").Copy()
$swiss.Range("A19").PasteSpecial(-4122)
$swiss.Range("A20").Copy()
$swiss.Range("A21").PasteSpecial(-4122)

# Market label
$swiss.Range("B2").Value = "Switzerland Market"

# Repeater codes for the Swiss market
$swiss.Range("A16").Value = "P32AR-CH"
$swiss.Range("A17").Value = "P32DR-CH"
$swiss.Range("A19").Value = "PR1DSCH"
$swiss.Range("A21").Value = "PR8ASCH"

# User story / NGC reference
$swiss.Range("B4").Value = "NGC-3476/T2345"

# Column widths for the new sheet
$swiss.Columns.Item(2).ColumnWidth = 22.833333333333332
$swiss.Columns.Item(4).ColumnWidth = 29

# Update the selection on the Swiss sheet and make it the active tab
$swiss.Range("B16").Select()
$swiss.Activate()
